$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert a new row at position 11 (pushes existing rows 11-17 down to 12-18),
# copying formatting from the row above (xlFormatFromLeftOrAbove) so the new
# row picks up the same cell style used by its neighbors, to add the new
# "charge_matrix" field between "peak_matrix_sem" and "spot_dist".
$ws.Rows(11).Insert(-4121, 0)

$ws.Range("A11").Value = "charge_matrix"
$ws.Range("B11").Value = "cell"
$ws.Range("C11").Value = "charge of each trace in trace_matrix_mean"

# Update selection to match the new active cell shown in the saved workbook.
$ws.Range("C11").Select()
